# Diag Classe un peu modifié et diag seq avancé
#
# Adds a short "next steps" note plus a class-diagram singleton remark right
# after the closing "L'As est une carte spéciale ..." paragraph (end of the
# card-rules section), and a trailing blank paragraph right after the
# existing bookmark paragraph, just before the section break.

$d = $word.ActiveDocument

# --- Locate the paragraph that closes the card-rules section -------------
# ("L'As est une carte spéciale qui permet ... le paquet de jeu. ")
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "L’As est une carte spéciale",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph ('L'As est une carte ...') not found"
}

$lastRuleIndex = $searchRange.Paragraphs(1).Index

# 1) Blank paragraph right after the card-rules section.
$lastRule = $d.Paragraphs($lastRuleIndex).Range
$lastRule.Collapse(0)
$lastRule.InsertParagraphAfter()

# 2) Short note about upcoming corrections (new paragraph).
$blank = $d.Paragraphs($lastRuleIndex + 1).Range
$blank.Collapse(0)
$blank.InsertParagraphAfter()
$notesPara = $d.Paragraphs($lastRuleIndex + 2)
$notesPara.Range.Text = "Nous savons déjà que certaine partie seront à corriger, comme le type de couleur, il y aura sûrement une définition de type avec une énumération. "

# 3) "DistribuerCarte : singleton" paragraph.
$notesEnd = $d.Paragraphs($lastRuleIndex + 2).Range
$notesEnd.Collapse(0)
$notesEnd.InsertParagraphAfter()
$singletonPara = $d.Paragraphs($lastRuleIndex + 3)
$singletonPara.Range.Text = "DistribuerCarte : singleton "

# --- Add a trailing blank paragraph right after the bookmark paragraph ---
# (the paragraph that only holds the _GoBack bookmark, just before the
# section break). Fall back to the document's last paragraph if, for some
# reason, that bookmark isn't present.
$bookmarkParaIndex = $d.Paragraphs.Count
for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
    if ($d.Bookmarks.Item($i).Name -eq "_GoBack") {
        $bookmarkParaIndex = $d.Bookmarks.Item($i).Range.Paragraphs(1).Index
        break
    }
}

$bookmarkPara = $d.Paragraphs($bookmarkParaIndex).Range
$bookmarkPara.Collapse(0)
$bookmarkPara.InsertParagraphAfter()

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
